# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# New columns AD:AF are appended after the existing data (which ends at AC),
# with the same header style as the rest of row 1, and every player row
# (2-42) gets the same team record: 74 wins, 88 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (AC1, which uses the
# bold/bordered/centered header style) onto the three new header cells so
# they pick up the same style instead of the default one.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Header labels for the new columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record values for every player row.
for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 74  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 88  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
